# Adding default_density as user input
# Insert two new parameter rows ("default_housing_density" and
# "default_vegetation_density") into the "parameters" sheet, just above the
# existing "fire_degradation_rate_min" row, and update
# fire_degradation_rate_min's value from 0.8 to 0.4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")

# Push the existing "fire_degradation_rate_min" (old row 15) and the rows
# below it down by two, freeing up rows 15 and 16 for the new parameters.
$ws.Rows.Item(15).Insert() | Out-Null
$ws.Rows.Item(16).Insert() | Out-Null

$ws.Range("A15").Value = "default_housing_density"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "required for case instance; in units/sq km, 1: 0; 2: <6; 3: 6-50; 4: 50-741; 5: > 741"

$ws.Range("A16").Value = "default_vegetation_density"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "required for case instance; 0 if <50% vegetated, 1 if >50% vegetated"

# fire_degradation_rate_min now lives on row 17; update its value.
$ws.Range("B17").Value = 0.4

# Match the author's final cursor position on the sheet.
$ws.Range("C17").Select() | Out-Null
